$p = $ppt.ActivePresentation
$s = $p.Slides.Item(13)

# 1) Split the "Seguimiento archivos y lectura contadoras - pi/2 BPSK, QPSK"
#    textbox into three runs, changing "archivos" -> "activos":
#      "Seguimiento " + "activos " + "y lectura contadoras - pi/2 BPSK, QPSK"
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq "Seguimiento archivos y lectura contadoras - pi/2 BPSK, QPSK") {
        $tr = $shp.TextFrame.TextRange
        $tr.Text = "Seguimiento "
        [void]$tr.InsertAfter("activos ")
        [void]$tr.InsertAfter("y lectura contadoras - pi/2 BPSK, QPSK")
        break
    }
}

# 2) Remove two stray/duplicate textboxes left over from the Google Slides
#    import: the lone "6." numbering box and the duplicated
#    "Sistemas buscapersonas - FSK" label (identified by their unique
#    shape Ids; iterate back-to-front so deleting doesn't disturb indices).
$idsToRemove = @(1108541324, 995512332)
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shp = $s.Shapes.Item($i)
    if ($idsToRemove -contains $shp.Id) {
        $shp.Delete()
    }
}
